# Apply the "demoData.xlsx" restructuring:
#   Sheet1   -> Create      (TC1000 create-user row + TC1002 create-user-no-name row)
#   testData -> List users  (TC1001 list-users row)
#   (new)    -> Update      (TC1003 update-user row)

$wb = $excel.ActiveWorkbook

$wsCreate = $wb.Worksheets.Item("Sheet1")
$wsList   = $wb.Worksheets.Item("testData")

# --- rename existing sheets -------------------------------------------------
$wsCreate.Name = "Create"
$wsList.Name   = "List users"

# --- add the new "Update" sheet after "List users" --------------------------
$wsUpdate = $wb.Worksheets.Add($null, $wsList)
$wsUpdate.Name = "Update"

# =============================================================================
# Create sheet: drop the endpoint/method columns in favour of name/job,
# remove the hyperlink, and add a second test-case row.
# =============================================================================
$wsCreate.Hyperlinks.Delete()

$wsCreate.Range("B1").Value = "name"
$wsCreate.Range("C1").Value = "job"

$wsCreate.Range("B2").Value = "test"
$wsCreate.Range("B2").Style = "Normal"
$wsCreate.Range("C2").Value = "test"

$wsCreate.Range("A3").Value = "TC1002"
$wsCreate.Range("C3").Value = "test"

$wsCreate.Range("L13").Select()

# =============================================================================
# List users sheet: only the testcase column remains, with a styled 2nd row.
# =============================================================================
$wsList.Range("A2").Copy()
$wsList.Range("B2").PasteSpecial(-4122)
$wsList.Application.CutCopyMode = $false

$wsList.Range("B1").ClearContents()
$wsList.Range("A2").Value = "TC1001"
$wsList.Range("B2").ClearContents()

$wsList.Range("A1:A2").Select()

# =============================================================================
# Update sheet: same header shape as Create, single "updated" data row.
# =============================================================================
$wsUpdate.Range("A1").Value = "testcase"
$wsUpdate.Range("B1").Value = "name"
$wsUpdate.Range("C1").Value = "job"

$wsCreate.Range("A2").Copy()
$wsUpdate.Range("A2").PasteSpecial(-4122)
$wsUpdate.Application.CutCopyMode = $false

$wsUpdate.Range("A2").Value = "TC1003"
$wsUpdate.Range("B2").Value = "updated"
$wsUpdate.Range("C2").Value = "updated"

$wsUpdate.Range("B2:C2").Select()
